$d = $word.ActiveDocument

function New-RunsXml {
    param([string]$InnerXml)
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $InnerXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# The "Tools" bullet list originally reads (all sub-bullets at ilvl 1):
#   Execute sql / Enrich city data / Explain query / Generate insights
# It becomes a two-level list where each original sub-bullet's text slides
# down into the next slot and picks up a new detail bullet underneath it,
# e.g. "Enrich city data" -> "Returns query, explanation, results, and row
# count" (demoted under "Execute sql"), and so on, ending with a brand new
# "Generate insights" bullet plus its own new detail line.

# ---------------------------------------------------------------------------
# Phase 1: reshape the paragraph list, first establishing every slot (blank
# placeholders for the brand-new paragraphs) while the indices are still
# predictable, before any content/XML is poured in.
# ---------------------------------------------------------------------------

# Paragraph 8 ("Enrich city data", ilvl 1) -> becomes the new detail bullet
# for "Execute sql": ilvl 2, "Returns query, explanation, results, and row
# count". Edited in place (keeps its original paragraph identity).
$p8 = $d.Paragraphs.Item(8)
$p8.Range.ListFormat.ListLevelNumber = 3
$p8.Range.Text = "Returns query, explanation, results, and row count"

# Paragraph 9 ("Explain query", ilvl 1) -> becomes "Enrich city data"
# (still ilvl 1). Edited in place.
$p9 = $d.Paragraphs.Item(9)
$p9.Range.Text = "Enrich city data"

# New blank paragraph 10 right after it (detail bullet slot for "Enrich
# city data").
$null = $p9.Range.InsertParagraphAfter()

# Paragraph 11 ("Generate insights", ilvl 1) -> becomes "Explain query"
# (still ilvl 1). Edited in place.
$p11 = $d.Paragraphs.Item(11)
$p11.Range.Text = "Explain query"

# New blank paragraph 12 right after it (detail bullet slot for "Explain
# query").
$null = $p11.Range.InsertParagraphAfter()

# New blank paragraph 13 (brand new top-level "Generate insights" slot).
$p12 = $d.Paragraphs.Item(12)
$null = $p12.Range.InsertParagraphAfter()

# New blank paragraph 14 (detail bullet slot for the new "Generate
# insights").
$p13 = $d.Paragraphs.Item(13)
$null = $p13.Range.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# Phase 2: pour content into the blank slots, working from the bottom of the
# document upward. Filling bottom-up guarantees every InsertXML target
# always has a following paragraph at the moment it runs (InsertXML behaves
# oddly -- leaves a stray empty paragraph behind -- when its target range is
# the very last paragraph in the document).
# ---------------------------------------------------------------------------

# Paragraph 14 (ilvl 2): "Contextual insights and follow up about
# analysis" -- plain text, no flagged words, so Range.Text is enough.
$p14 = $d.Paragraphs.Item(14)
$p14.Range.ListFormat.ListLevelNumber = 3
$p14.Range.Text = "Contextual insights and follow up about analysis"

# Paragraph 13 (ilvl 1): "Generate insights" -- plain text too.
$p13 = $d.Paragraphs.Item(13)
$p13.Range.ListFormat.ListLevelNumber = 2
$p13.Range.Text = "Generate insights"

# Paragraph 12 (ilvl 2): "Explain query and sql to users who don't know
# sql" -- both "sql" occurrences are flagged by the spell checker, so build
# it with explicit runs/proofErr markers via InsertXML.
$p12 = $d.Paragraphs.Item(12)
$p12.Range.ListFormat.ListLevelNumber = 3
$xml12 = New-RunsXml('<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Explain query and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sql</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to users who don&#8217;t know </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sql</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body>')
$null = $p12.Range.InsertXML($xml12)

# Paragraph 10 (ilvl 2): "Weather, timezone, demographics, economics" --
# "timezone" is flagged by the spell checker.
$p10 = $d.Paragraphs.Item(10)
$p10.Range.ListFormat.ListLevelNumber = 3
$xml10 = New-RunsXml('<w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Weather, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>timezone</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, demographics, economics</w:t></w:r></w:p></w:body>')
$null = $p10.Range.InsertXML($xml10)
